$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''58.139.07'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '''3.113.72'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''527.95'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = '''141.70'
$ws.Range('E6').Value = '  -1.73%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '''3.113.73'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = '''0.445'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').Value = '''7.18'
$ws.Range('E10').Value = '  -2.96%  '
$ws.Range('D11').Value = '''0.109'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').Value = '''3.646.21'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').Value = '''25.73'
$ws.Range('E15').Value = '  -4.92%  '
$ws.Range('D16').Value = '''0.0000166'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '''58.169.24'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '''3.098.69'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = '''12.77'
$ws.Range('E20').Value = '  -2.07%  '
$ws.Range('D21').Value = '''8.00'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').Value = '''343.21'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '''0.515'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').Value = '''67.71'
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '''0.0₃0932'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '''6.38'
$ws.Range('E30').Value = '  -5.83%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '''7.29'
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').Value = '''1.88'
$ws.Range('E32').Value = '  +1.46%  '
$ws.Range('D33').Value = '''21.11'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('D35').Value = '''158.48'
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('D36').Value = '''4.64'
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').Value = '''6.18'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '''26.31'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').Value = '''1.24'
$ws.Range('E39').Value = '  -4.25%  '
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('E41').Value = '  +10.67%  '
$ws.Range('D42').Value = '''4.01'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').Value = '''0.690'
$ws.Range('E43').Value = '  +3.71%  '
$ws.Range('D44').Value = '''3.152.90'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = '''0.0263'
$ws.Range('E47').Value = '  +2.25%  '
$ws.Range('D48').Value = '''2.269.88'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('D49').Value = '''0.997'
$ws.Range('E49').Value = '  +2.98%  '
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').Value = '''20.60'
$ws.Range('E51').Value = '  -1.31%  '
